$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.300.59"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3
$ws.Range("D3").Value = "1.867.93"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.91%  "

# Row 6
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4725"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.65%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06480"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07734"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.80%  "

# Row 13
$ws.Range("D13").Value = "1.870.03"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7036"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.101"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "275.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18
$ws.Range("D18").Value = "30.280.59"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007561"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").Value = "2.117.77"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.207"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.136"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.311"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.909"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "

# Row 30
$ws.Range("E30").Value = "  +1.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09863"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.250"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.72%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.031"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04736"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6916"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.706"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "

# Row 39
$ws.Range("E39").Value = "  -1.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.743"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.337"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.76%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.900"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.44%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4087"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.281"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.068"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.68%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "922.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
